$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.540.84"
$ws.Range("E2").Value = "  +1.86%  "

$ws.Range("D3").Value = "2.551.37"
$ws.Range("E3").Value = "  +5.08%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.14%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "571.49"
$ws.Range("E5").Value = "  +2.65%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "150.76"
$ws.Range("E6").Value = "  +8.82%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.999"
$ws.Range("E7").Value = "  -0.11%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.587"
$ws.Range("E8").Value = "  +0.32%  "

$ws.Range("D9").Value = "2.545.39"
$ws.Range("E9").Value = "  +4.90%  "

$ws.Range("E10").Value = "  +2.29%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.75"
$ws.Range("E11").Value = "  +0.10%  "

$ws.Range("E12").Value = "  +1.05%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.359"
$ws.Range("E13").Value = "  +3.27%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "28.46"
$ws.Range("E14").Value = "  +8.33%  "

$ws.Range("D15").Value = "3.001.66"
$ws.Range("E15").Value = "  +4.81%  "

$ws.Range("D16").Value = "63.440.75"
$ws.Range("E16").Value = "  +1.83%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.0000145"
$ws.Range("E17").Value = "  +2.82%  "

$ws.Range("D18").Value = "2.567.31"
$ws.Range("E18").Value = "  +5.74%  "

$ws.Range("E19").Value = "  +4.55%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "340.98"
$ws.Range("E20").Value = "  -1.49%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.38"
$ws.Range("E21").Value = "  +4.37%  "

$ws.Range("E22").Value = "  +0.74%  "

$ws.Range("E23").Value = "  +0.15%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "66.19"
$ws.Range("E24").Value = "  +1.64%  "

$ws.Range("E25").Value = "  -0.96%  "

$ws.Range("E26").Value = "  +5.85%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.63"
$ws.Range("E27").Value = "  +5.77%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.51"
$ws.Range("E28").Value = "  +12.15%  "

$ws.Range("E29").Value = "  +0.03%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.23"
$ws.Range("E30").Value = "  +12.81%  "

$ws.Range("D31").Value = "0.0₃0831"
$ws.Range("E31").Value = "  +6.17%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.88"
$ws.Range("E32").Value = "  +4.31%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "178.16"
$ws.Range("E33").Value = "  +3.71%  "

$ws.Range("E34").Value = "  +8.32%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "420.83"
$ws.Range("E35").Value = "  +10.63%  "

$ws.Range("E36").Value = "  +2.32%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "19.18"
$ws.Range("E37").Value = "  +3.38%  "

$ws.Range("E38").Value = "  +0.22%  "

$ws.Range("E40").Value = "  +5.88%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.999"
$ws.Range("E41").Value = "  -0.11%  "

$ws.Range("E42").Value = "  +1.44%  "

$ws.Range("E43").Value = "  +6.40%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "3.83"
$ws.Range("E44").Value = "  +4.71%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "21.13"
$ws.Range("E45").Value = "  +2.21%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.612"
$ws.Range("E46").Value = "  +3.89%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0530"
$ws.Range("E47").Value = "  +2.22%  "

$ws.Range("E48").Value = "  +1.76%  "

$ws.Range("E49").Value = "  +7.80%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "18.70"
$ws.Range("E50").Value = "  +4.58%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.83"
$ws.Range("E51").Value = "  +8.24%  "
